$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1269.2727
$ws.Range("I40").Value = 1120
$ws.Range("J40").Value = 1667.3334
$ws.Range("K40").Value = 1120
$ws.Range("L40").Value = 1667.3334
$ws.Range("M40").Value = -945
$ws.Range("N40").Value = -2017.3334
$ws.Range("H43").Value = 696.7059
$ws.Range("I43").Value = 773.61536
$ws.Range("J43").Value = 649.0952
$ws.Range("K43").Value = 773.61536
$ws.Range("L43").Value = 649.0952
$ws.Range("M43").Value = -704.61536
$ws.Range("N43").Value = -787.0952
$ws.Range("H113").Value = 2868.5356
$ws.Range("I113").Value = 2305.6667
$ws.Range("J113").Value = 3518
$ws.Range("K113").Value = 2305.6667
$ws.Range("L113").Value = 3518
$ws.Range("M113").Value = 948.3332999999998
$ws.Range("N113").Value = -10026
$ws.Range("H133").Value = 49975
$ws.Range("J133").Value = 49975
$ws.Range("L133").Value = 49975
$ws.Range("N133").Value = -60095
$ws.Range("H137").Value = 901.29266
$ws.Range("I137").Value = 726.3913
$ws.Range("K137").Value = 2179.1739
$ws.Range("M137").Value = 370.8261000000002
$ws.Range("H138").Value = 2455.125
$ws.Range("I138").Value = 2027.7715
$ws.Range("J138").Value = 2859.3784
$ws.Range("K138").Value = 6083.3145
$ws.Range("L138").Value = 8578.135200000001
$ws.Range("M138").Value = -943.3145000000004
$ws.Range("N138").Value = -18858.1352
$ws.Range("H139").Value = 63786.668
$ws.Range("J139").Value = 63786.668
$ws.Range("L139").Value = 63786.668
$ws.Range("N139").Value = -74066.66800000001
$ws.Range("H140").Value = 84189.09
$ws.Range("J140").Value = 84189.09
$ws.Range("L140").Value = 84189.09
$ws.Range("N140").Value = -94549.09

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 391510.25
$ws.Range("I32").Value = 4560.631
$ws.Range("J32").Value = 3186146.5
$ws.Range("K32").Value = 4560.631
$ws.Range("L32").Value = 3186146.5
$ws.Range("M32").Value = -4273.631
$ws.Range("N32").Value = -3186720.5
$ws.Range("H74").Value = 2593
$ws.Range("I74").Value = 2380.6875
$ws.Range("J74").Value = 2932.7
$ws.Range("K74").Value = 2380.6875
$ws.Range("L74").Value = 2932.7
$ws.Range("M74").Value = -1506.6875
$ws.Range("N74").Value = -4680.7
$ws.Range("H77").Value = 2593
$ws.Range("I77").Value = 2380.6875
$ws.Range("J77").Value = 2932.7
$ws.Range("K77").Value = 11903.4375
$ws.Range("L77").Value = 14663.5
$ws.Range("M77").Value = -7535.4375
$ws.Range("N77").Value = -23399.5
$ws.Range("H132").Value = 2623.1667
$ws.Range("I132").Value = 2139.5881
$ws.Range("J132").Value = 3255.5386
$ws.Range("K132").Value = 6418.7643
$ws.Range("L132").Value = 9766.6158
$ws.Range("M132").Value = -3888.7643
$ws.Range("N132").Value = -14826.6158
$ws.Range("H138").Value = 60825.715
$ws.Range("J138").Value = 60825.715
$ws.Range("L138").Value = 60825.715
$ws.Range("N138").Value = -71105.715
$ws.Range("H139").Value = 64500
$ws.Range("J139").Value = 64500
$ws.Range("L139").Value = 64500
$ws.Range("N139").Value = -74780

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 48542.855
$ws.Range("J81").Value = 48542.855
$ws.Range("L81").Value = 48542.855
$ws.Range("N81").Value = -50664.855
$ws.Range("H84").Value = 48542.855
$ws.Range("J84").Value = 48542.855
$ws.Range("L84").Value = 145628.565
$ws.Range("N84").Value = -156236.565
$ws.Range("H94").Value = 1586.6875
$ws.Range("I94").Value = 1561.5714
$ws.Range("J94").Value = 1606.2222
$ws.Range("K94").Value = 1561.5714
$ws.Range("L94").Value = 1606.2222
$ws.Range("M94").Value = -1110.5714
$ws.Range("N94").Value = -2508.2222
$ws.Range("H132").Value = 50755
$ws.Range("J132").Value = 50755
$ws.Range("L132").Value = 50755
$ws.Range("N132").Value = -60875
$ws.Range("H138").Value = 47651.43
$ws.Range("J138").Value = 47651.43
$ws.Range("L138").Value = 47651.43
$ws.Range("N138").Value = -57931.43
$ws.Range("H140").Value = 73200
$ws.Range("J140").Value = 73200
$ws.Range("L140").Value = 73200
$ws.Range("N140").Value = -83560

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1196.3636
$ws.Range("J58").Value = 1417.8667
$ws.Range("L58").Value = 1417.8667
$ws.Range("N58").Value = -1823.8667
$ws.Range("H136").Value = 1196.3636
$ws.Range("J136").Value = 1417.8667
$ws.Range("L136").Value = 4253.6001
$ws.Range("N136").Value = -9353.6001
$ws.Range("H140").Value = 89966.664
$ws.Range("J140").Value = 89966.664
$ws.Range("L140").Value = 89966.664
$ws.Range("N140").Value = -100326.664
$ws.Range("H141").Value = 38228.57
$ws.Range("J141").Value = 36300
$ws.Range("L141").Value = 36300
$ws.Range("N141").Value = -46660

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 500.5
$ws.Range("I107").Value = 292
$ws.Range("J107").Value = 848
$ws.Range("K107").Value = 876
$ws.Range("L107").Value = 2544
$ws.Range("M107").Value = 1044
$ws.Range("N107").Value = -6384
$ws.Range("H122").Value = 571.7917
$ws.Range("J122").Value = 918.6667
$ws.Range("L122").Value = 8268.0003
$ws.Range("N122").Value = -13168.0003
$ws.Range("H124").Value = 5000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 5000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 15000
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -24820
$ws.Range("H131").Value = 834.76
$ws.Range("I131").Value = 450
$ws.Range("J131").Value = 842.61224
$ws.Range("K131").Value = 1350
$ws.Range("L131").Value = 2527.83672
$ws.Range("M131").Value = 3690
$ws.Range("N131").Value = -12607.83672
$ws.Range("H132").Value = 921.08
$ws.Range("I132").Value = 1004.5455
$ws.Range("J132").Value = 855.5
$ws.Range("K132").Value = 9040.9095
$ws.Range("L132").Value = 7699.5
$ws.Range("M132").Value = -6510.9095
$ws.Range("N132").Value = -12759.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 34466
$ws.Range("J133").Value = 34466
$ws.Range("L133").Value = 34466
$ws.Range("N133").Value = -44586
$ws.Range("H138").Value = 61708.332
$ws.Range("J138").Value = 61708.332
$ws.Range("L138").Value = 61708.332
$ws.Range("N138").Value = -71988.33199999999
$ws.Range("H140").Value = 99733
$ws.Range("J140").Value = 99733
$ws.Range("L140").Value = 99733
$ws.Range("N140").Value = -110093
$ws.Range("H141").Value = 66577.25
$ws.Range("J141").Value = 66577.25
$ws.Range("L141").Value = 66577.25
$ws.Range("N141").Value = -76937.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1903
$ws.Range("I93").Value = 1356.381
$ws.Range("J93").Value = 2786
$ws.Range("K93").Value = 1356.381
$ws.Range("L93").Value = 2786
$ws.Range("M93").Value = -108.3810000000001
$ws.Range("N93").Value = -5282
$ws.Range("H132").Value = 2693.9836
$ws.Range("I132").Value = 2833.28
$ws.Range("J132").Value = 2060.818
$ws.Range("K132").Value = 8499.84
$ws.Range("L132").Value = 6182.454000000001
$ws.Range("M132").Value = -5969.84
$ws.Range("N132").Value = -11242.454
$ws.Range("H138").Value = 64679.57
$ws.Range("J138").Value = 64679.57
$ws.Range("L138").Value = 64679.57
$ws.Range("N138").Value = -74959.57000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 30303610
$ws.Range("I113").Value = 489.9643
$ws.Range("J113").Value = 200001090
$ws.Range("K113").Value = 1469.8929
$ws.Range("L113").Value = 600003270
$ws.Range("M113").Value = 700.1071000000002
$ws.Range("N113").Value = -600007610
$ws.Range("H141").Value = 54470
$ws.Range("J141").Value = 55188.89
$ws.Range("L141").Value = 55188.89
$ws.Range("N141").Value = -65548.89
